$d = $word.ActiveDocument
$c = $d.Content
$found = $c.Find.Execute("Le bénévole ou l’employé valide", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "found=$found start=$($c.Start) end=$($c.End) text=[$($c.Text)]"

$len = ("Le bénévole ou l’employé ").Length
$r2 = $d.Range($c.Start, $c.Start + $len)
Write-Output "r2 text=[$($r2.Text)]"

$found2 = $r2.Find.Execute("Le bénévole ou l’employé ", $true, $false, $false, $false, $false, $true, 1, $false, "Le bénévole ou l’employé ", 1)
Write-Output "found2=$found2"
Write-Output "doc around=[$($d.Range($c.Start-5, $c.Start+30).Text)]"
